# ============================================================
# Update report header text (Volume/Number and week-covering dates)
# ============================================================
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Volume 30   Number  37" -> "...38"  (A8, shared string with 4 runs)
$volCell = $ws.Range("A8")
$volCell.Characters(21, 2).Text = "38"

# "Report Covering the Week  9/11/2023  Through  9/17/2023" (C9)
$weekCell = $ws.Range("C9")
$weekCell.Characters(27, 9).Text = "9/18/2023"
$weekCell.Characters(47, 9).Text = "9/24/2023"

# ============================================================
# Update the 42nd Precinct weekly crime-complaint data table
# (rows 14-29, columns C:N) with newly collected figures
# ============================================================

$ws.Range("L14").Value = -47.058823529411
$ws.Range("N14").Value = -60.869565217391
$ws.Range("C14").Copy($ws.Range("D15"))
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("C14").Copy($ws.Range("F15"))
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = -100
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -42.857142857142
$ws.Range("F16").Value = 28
$ws.Range("G16").Value = 34
$ws.Range("H16").Value = -17.647058823529
$ws.Range("I16").Value = 305
$ws.Range("J16").Value = 321
$ws.Range("K16").Value = -4.984423676012
$ws.Range("L16").Value = 31.465517241379
$ws.Range("M16").Value = 47.342995169082
$ws.Range("N16").Value = -66.298342541436
$ws.Range("C17").Value = 14
$ws.Range("D17").Value = 16
$ws.Range("E17").Value = -12.5
$ws.Range("F17").Value = 60
$ws.Range("G17").Value = 63
$ws.Range("H17").Value = -4.761904761904
$ws.Range("I17").Value = 530
$ws.Range("J17").Value = 497
$ws.Range("K17").Value = 6.639839034205
$ws.Range("L17").Value = 13.733905579399
$ws.Range("M17").Value = 130.434782608696
$ws.Range("N17").Value = -27.397260273972
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 27
$ws.Range("H18").Value = -44.444444444444
$ws.Range("I18").Value = 145
$ws.Range("J18").Value = 251
$ws.Range("K18").Value = -42.231075697211
$ws.Range("L18").Value = 2.112676056338
$ws.Range("M18").Value = 30.630630630630
$ws.Range("N18").Value = -81.691919191919
$ws.Range("C19").Value = 17
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 70
$ws.Range("F19").Value = 55
$ws.Range("G19").Value = 44
$ws.Range("H19").Value = 25
$ws.Range("I19").Value = 376
$ws.Range("J19").Value = 374
$ws.Range("K19").Value = 0.534759358288
$ws.Range("L19").Value = 14.634146341463
$ws.Range("M19").Value = 111.23595505618
$ws.Range("N19").Value = 41.353383458646
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 75
$ws.Range("F20").Value = 31
$ws.Range("G20").Value = 27
$ws.Range("H20").Value = 14.814814814814
$ws.Range("I20").Value = 336
$ws.Range("J20").Value = 223
$ws.Range("K20").Value = 50.672645739910
$ws.Range("L20").Value = 121.052631578947
$ws.Range("M20").Value = 314.814814814815
$ws.Range("N20").Value = -6.666666666666
$ws.Range("C21").Value = 45
$ws.Range("D21").Value = 43
$ws.Range("E21").Value = 4.651162790697
$ws.Range("F21").Value = 189
$ws.Range("G21").Value = 196
$ws.Range("H21").Value = -3.571428571428
$ws.Range("I21").Value = 1724
$ws.Range("J21").Value = 1701
$ws.Range("K21").Value = 1.352145796590
$ws.Range("L21").Value = 26.485693323551
$ws.Range("M21").Value = 108.464328899637
$ws.Range("N21").Value = -45.007974481658
$ws.Range("C14").Copy($ws.Range("C22"))
$ws.Range("C23").Value = 8
$ws.Range("D23").Value = 9
$ws.Range("E23").Value = -11.111111111111
$ws.Range("F23").Value = 33
$ws.Range("G23").Value = 34
$ws.Range("H23").Value = -2.941176470588
$ws.Range("I23").Value = 305
$ws.Range("J23").Value = 263
$ws.Range("K23").Value = 15.969581749049
$ws.Range("L23").Value = 103.333333333333
$ws.Range("M23").Value = 114.788732394366
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 33
$ws.Range("E24").Value = -6.060606060606
$ws.Range("G24").Value = 118
$ws.Range("H24").Value = -19.491525423728
$ws.Range("I24").Value = 877
$ws.Range("J24").Value = 950
$ws.Range("K24").Value = -7.684210526315
$ws.Range("L24").Value = 27.285921625544
$ws.Range("M24").Value = 57.450628366247
$ws.Range("C25").Value = 19
$ws.Range("D25").Value = 22
$ws.Range("E25").Value = -13.636363636363
$ws.Range("F25").Value = 75
$ws.Range("G25").Value = 74
$ws.Range("H25").Value = 1.351351351351
$ws.Range("I25").Value = 805
$ws.Range("J25").Value = 738
$ws.Range("K25").Value = 9.078590785907
$ws.Range("L25").Value = 24.613003095975
$ws.Range("M25").Value = 26.771653543307
$ws.Range("I14").Copy($ws.Range("C26"))
$ws.Range("C26").Value = 2
$ws.Range("C14").Copy($ws.Range("D26"))
$ws.Range("E14").Copy($ws.Range("E26"))
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 200
$ws.Range("I26").Value = 41
$ws.Range("K26").Value = -10.869565217391
$ws.Range("L26").Value = -2.380952380952
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = 50
$ws.Range("F27").Value = 13
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 333.333333333333
$ws.Range("I27").Value = 74
$ws.Range("J27").Value = 56
$ws.Range("K27").Value = 32.142857142857
$ws.Range("L27").Value = 57.446808510638
$ws.Range("I14").Copy($ws.Range("C28"))
$ws.Range("C28").Value = 2
$ws.Range("C14").Copy($ws.Range("D28"))
$ws.Range("E14").Copy($ws.Range("E28"))
$ws.Range("F28").Value = 4
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 33
$ws.Range("K28").Value = 6.451612903225
$ws.Range("L28").Value = -45.901639344262
$ws.Range("M28").Value = -10.810810810810
$ws.Range("N28").Value = -58.227848101265
$ws.Range("I14").Copy($ws.Range("C29"))
$ws.Range("C29").Value = 1
$ws.Range("C14").Copy($ws.Range("D29"))
$ws.Range("E14").Copy($ws.Range("E29"))
$ws.Range("F29").Value = 3
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 27
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = -47.058823529411
$ws.Range("M29").Value = -12.903225806451
$ws.Range("N29").Value = -65.822784810126
